$d = $word.ActiveDocument

# 1. Update the date in the letter header.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address paragraph ("130 Baroni Ave., San Jose CA 95136")
#    into two paragraphs: "130 Baroni Ave." and "San Jose, CA 95136".
#    Locate the specific paragraph right after "Daniel Guillen Jr" (the addressee
#    block), not the similar text inside the PROPERTY ADDRESS table.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "130 Baroni Ave., San Jose CA 95136`r") {
        $rng = $para.Range
        $rng.Find.Execute(", San Jose CA 95136", $false, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
        $rng.Text = "`r" + "San Jose, CA 95136"

        $newPara = $d.Paragraphs.Item($i + 1)
        $f = $newPara.Range.Font
        $f.Name = "Arial"
        $f.Size = 11
        $f.NameBi = "Arial"
        $f.SizeBi = 11
        break
    }
}

# 3. Remove the now-redundant blank "No Spacing" paragraph that used to sit
#    directly below "... Board of Directors".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Baroni Green Homeowners Association Board of Directors`r") {
        $blank = $d.Paragraphs.Item($i + 1)
        if ($blank.Range.Text -eq "`r") {
            $blank.Range.Delete()
        }
        break
    }
}
